$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 25923
$ws.Range("E2").Value = 505056033865
$ws.Range("F2").Value = 4418612806
$ws.Range("G2").Value = 0.40406
$ws.Range("D3").Value = 1633.2
$ws.Range("E3").Value = 196426979119
$ws.Range("F3").Value = 3549291287
$ws.Range("G3").Value = -0.03559
$ws.Range("D4").Value = 0.9994
$ws.Range("E4").Value = 82891856881
$ws.Range("F4").Value = 7954811634
$ws.Range("G4").Value = -0.01768
$ws.Range("D5").Value = 214.48
$ws.Range("E5").Value = 33014771455
$ws.Range("F5").Value = 229554722
$ws.Range("G5").Value = -0.06374
$ws.Range("D6").Value = 0.502923
$ws.Range("E6").Value = 26645139059
$ws.Range("F6").Value = 528038859
$ws.Range("G6").Value = 1.05653
$ws.Range("D7").Value = 0.999471
$ws.Range("E7").Value = 26183102929
$ws.Range("F7").Value = 1975122872
$ws.Range("G7").Value = -0.05115
$ws.Range("D8").Value = 1632.61
$ws.Range("E8").Value = 13951987479
$ws.Range("F8").Value = 6404663
$ws.Range("G8").Value = -0.02969
$ws.Range("D9").Value = 0.254568
$ws.Range("E9").Value = 8925549096
$ws.Range("F9").Value = 108085140
$ws.Range("G9").Value = -0.2249
$ws.Range("D10").Value = 0.063301
$ws.Range("E10").Value = 8918233250
$ws.Range("F10").Value = 222858490
$ws.Range("G10").Value = -0.09375
$ws.Range("D11").Value = 19.4
$ws.Range("E11").Value = 7939475440
$ws.Range("F11").Value = 195136171
$ws.Range("G11").Value = -0.84099
$ws.Range("D12").Value = 0.076975
$ws.Range("E12").Value = 6877035415
$ws.Range("F12").Value = 137841026
$ws.Range("G12").Value = 0.45218
$ws.Range("B13").Value = "TON"
$ws.Range("C13").Value = "Toncoin"
$ws.Range("D13").Value = 1.87
$ws.Range("E13").Value = 6427508543
$ws.Range("F13").Value = 28787534
$ws.Range("G13").Value = -2.41439
$ws.Range("B14").Value = "DOT"
$ws.Range("C14").Value = "Polkadot"
$ws.Range("D14").Value = 4.27
$ws.Range("E14").Value = 5430187963
$ws.Range("F14").Value = 68667287
$ws.Range("G14").Value = 0.81496
$ws.Range("D15").Value = 0.543806
$ws.Range("E15").Value = 5071450465
$ws.Range("F15").Value = 132458424
$ws.Range("G15").Value = 0.08871
$ws.Range("D16").Value = 64.06999999999999
$ws.Range("E16").Value = 4717511503
$ws.Range("F16").Value = 258219481
$ws.Range("G16").Value = 0.15161
$ws.Range("D17").Value = [double]"7.8e-06"
$ws.Range("E17").Value = 4603976979
$ws.Range("F17").Value = 66639681
$ws.Range("G17").Value = -0.72229
$ws.Range("D18").Value = 25927
$ws.Range("E18").Value = 4224992976
$ws.Range("F18").Value = 25511543
$ws.Range("G18").Value = 0.34156
$ws.Range("D19").Value = 0.999483
$ws.Range("E19").Value = 3874000837
$ws.Range("F19").Value = 46237840
$ws.Range("G19").Value = -0.16045
$ws.Range("D20").Value = 195.13
$ws.Range("E20").Value = 3806659582
$ws.Range("F20").Value = 119764334
$ws.Range("G20").Value = -2.42734
$ws.Range("B21").Value = "LEO"
$ws.Range("C21").Value = "LEO Token"
$ws.Range("D21").Value = 3.9
$ws.Range("E21").Value = 3629732023
$ws.Range("F21").Value = 138620
$ws.Range("G21").Value = -0.01904
$ws.Range("B22").Value = "AVAX"
$ws.Range("C22").Value = "Avalanche"
$ws.Range("D22").Value = 9.92
$ws.Range("E22").Value = 3512172570
$ws.Range("F22").Value = 95226789
$ws.Range("G22").Value = 0.7204199999999999
$ws.Range("B23").Value = "TUSD"
$ws.Range("C23").Value = "TrueUSD"
$ws.Range("D23").Value = 0.998293
$ws.Range("E23").Value = 3443012753
$ws.Range("F23").Value = 878740664
$ws.Range("G23").Value = 0.00778
$ws.Range("B24").Value = "UNI"
$ws.Range("C24").Value = "Uniswap"
$ws.Range("D24").Value = 4.44
$ws.Range("E24").Value = 3348371247
$ws.Range("F24").Value = 61143801
$ws.Range("G24").Value = 1.8923
$ws.Range("D25").Value = 6.03
$ws.Range("E25").Value = 3246988179
$ws.Range("F25").Value = 176635621
$ws.Range("G25").Value = 1.48924
$ws.Range("B26").Value = "XLM"
$ws.Range("C26").Value = "Stellar"
$ws.Range("D26").Value = 0.114821
$ws.Range("E26").Value = 3155281863
$ws.Range("F26").Value = 59611121
$ws.Range("G26").Value = 1.98807
$ws.Range("B27").Value = "BUSD"
$ws.Range("C27").Value = "Binance USD"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 2857631613
$ws.Range("F27").Value = 1401855790
$ws.Range("G27").Value = -0.01037
$ws.Range("D28").Value = 140.69
$ws.Range("E28").Value = 2555100174
$ws.Range("F28").Value = 46399656
$ws.Range("G28").Value = 0.13799
$ws.Range("D29").Value = 42.04
$ws.Range("E29").Value = 2524796778
$ws.Range("F29").Value = 5375862
$ws.Range("G29").Value = -1.60429
$ws.Range("D30").Value = 15.46
$ws.Range("E30").Value = 2208522245
$ws.Range("F30").Value = 41431556
$ws.Range("G30").Value = 0.31018
$ws.Range("D31").Value = 6.85
$ws.Range("E31").Value = 2004316497
$ws.Range("F31").Value = 78929214
$ws.Range("G31").Value = 2.01979
$ws.Range("D32").Value = 0.04984048
$ws.Range("E32").Value = 1656547219
$ws.Range("F32").Value = 52399637
$ws.Range("G32").Value = 1.29392
$ws.Range("B33").Value = "MNT"
$ws.Range("C33").Value = "Mantle"
$ws.Range("D33").Value = 0.452298
$ws.Range("E33").Value = 1462493425
$ws.Range("F33").Value = 7834196
$ws.Range("G33").Value = -0.02991
$ws.Range("D34").Value = 100.01
$ws.Range("E34").Value = 1454863824
$ws.Range("F34").Value = 9057600
$ws.Range("G34").Value = 0.37744
$ws.Range("B35").Value = "ICP"
$ws.Range("C35").Value = "Internet Computer"
$ws.Range("E35").Value = 1435096866
$ws.Range("F35").Value = 11741138
$ws.Range("G35").Value = -0.71478
$ws.Range("B36").Value = "FIL"
$ws.Range("C36").Value = "Filecoin"
$ws.Range("D36").Value = 3.18
$ws.Range("E36").Value = 1416937040
$ws.Range("F36").Value = 65446688
$ws.Range("G36").Value = 0.107
$ws.Range("D37").Value = 1.54
$ws.Range("E37").Value = 1366131615
$ws.Range("F37").Value = 33125198
$ws.Range("G37").Value = -0.18547
$ws.Range("D38").Value = 0.05039
$ws.Range("E38").Value = 1324134397
$ws.Range("F38").Value = 4280772
$ws.Range("G38").Value = 0.28518
$ws.Range("D39").Value = 5.55
$ws.Range("E39").Value = 1270569880
$ws.Range("F39").Value = 47496346
$ws.Range("G39").Value = -0.44056
$ws.Range("D40").Value = 0.893548
$ws.Range("E40").Value = 1140722630
$ws.Range("F40").Value = 72456501
$ws.Range("G40").Value = 0.60575
$ws.Range("D41").Value = 0.01546759
$ws.Range("E41").Value = 1126488234
$ws.Range("F41").Value = 29302542
$ws.Range("G41").Value = -0.38207
$ws.Range("B42").Value = "OP"
$ws.Range("C42").Value = "Optimism"
$ws.Range("D42").Value = 1.34
$ws.Range("E42").Value = 1067077806
$ws.Range("F42").Value = 64205586
$ws.Range("G42").Value = 0.38549
$ws.Range("B43").Value = "NEAR"
$ws.Range("C43").Value = "NEAR Protocol"
$ws.Range("D43").Value = 1.13
$ws.Range("E43").Value = 1064718009
$ws.Range("F43").Value = 38574968
$ws.Range("G43").Value = 0.21191
$ws.Range("D44").Value = 1129.01
$ws.Range("E44").Value = 1017292323
$ws.Range("F44").Value = 66882826
$ws.Range("G44").Value = -2.01109
$ws.Range("D45").Value = 1769.6
$ws.Range("E45").Value = 926134410
$ws.Range("F45").Value = 3521798
$ws.Range("G45").Value = -0.11296
$ws.Range("B46").Value = "XDC"
$ws.Range("C46").Value = "XDC Network"
$ws.Range("D46").Value = 0.058294
$ws.Range("E46").Value = 808335912
$ws.Range("F46").Value = 11862776
$ws.Range("G46").Value = 0.25537
$ws.Range("B47").Value = "FRAX"
$ws.Range("C47").Value = "Frax"
$ws.Range("D47").Value = 0.997702
$ws.Range("E47").Value = 805006687
$ws.Range("F47").Value = 4835283
$ws.Range("G47").Value = -0.08524
$ws.Range("B48").Value = "GRT"
$ws.Range("C48").Value = "The Graph"
$ws.Range("D48").Value = 0.087434
$ws.Range("E48").Value = 803609370
$ws.Range("F48").Value = 26240002
$ws.Range("G48").Value = 1.89629
$ws.Range("B49").Value = "AAVE"
$ws.Range("C49").Value = "Aave"
$ws.Range("D49").Value = 54.61
$ws.Range("E49").Value = 794272734
$ws.Range("F49").Value = 61599885
$ws.Range("G49").Value = 0.14874
$ws.Range("B50").Value = "WBT"
$ws.Range("C50").Value = "WhiteBIT Coin"
$ws.Range("D50").Value = 5.34
$ws.Range("E50").Value = 770440438
$ws.Range("F50").Value = 6008766
$ws.Range("G50").Value = -0.10526
$ws.Range("B51").Value = "ALGO"
$ws.Range("C51").Value = "Algorand"
$ws.Range("D51").Value = 0.094304
$ws.Range("E51").Value = 739219103
$ws.Range("F51").Value = 14335328
$ws.Range("G51").Value = 1.71761
